$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Fill in new "No internet" row (row 10)
$ws.Range("A10").Value = "No internet Implemented"
$ws.Range("B10").Value = 2
$ws.Range("C10").Value = 45500
$ws.Range("D10").Value = "When there is no internet app fetches houses from local storage, map shows warning."

# Match the row height / formatting used by the other populated rows (e.g. row 4/7)
$ws.Rows.Item(10).RowHeight = 30

$ws.Range("A4").Copy() | Out-Null
$ws.Range("A10").PasteSpecial(-4122) | Out-Null
$ws.Range("B7").Copy() | Out-Null
$ws.Range("B10").PasteSpecial(-4122) | Out-Null
$ws.Range("C7").Copy() | Out-Null
$ws.Range("C10").PasteSpecial(-4122) | Out-Null
$ws.Range("D7").Copy() | Out-Null
$ws.Range("D10").PasteSpecial(-4122) | Out-Null
$excel.CutCopyMode = 0

# Update view selection / scroll to match author's final state
$ws.Application.ActiveWindow.ScrollRow = 3
$ws.Range("A15").Select()

$wb.Save()
